# Status_v0102 fix note added for the two newest issues:
#   D-0020587 <EVT FW SSM CPLD> Most pattern of LED blinking don't match the
#       corresponding toggle rate.                         -> row 12
#   D-0020588 <EVT FW SSM CPLD> Mapping of Register 11h,13h and LEDs are
#       incorrect.                                          -> row 13
# Both rows get a new "Fix it in Status_v0102" note in column H, styled like
# the other "Fix it in ..." notes already present in the sheet (e.g. H8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new fix notes
$ws.Range("H12").Value2 = "Fix it in Status_v0102"
$ws.Range("H13").Value2 = "Fix it in Status_v0102"

# Match the formatting used by the existing "Fix it in ..." notes (red text,
# text number format) by copying the format from H8 onto the new cells.
$ws.Range("H8").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Reflect the scrolled-down view / active selection recorded for the sheet
# after the edit.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H15").Select() | Out-Null
